# Correct ES formality and adjust redcap instructions
#
# The intro paragraph before the Likert-scale table used informal
# "tu/te/ti" Spanish address forms and was split across many runs with
# redundant ". " separators. Re-word it to the formal "usted/le/le"
# register, merge the stray runs back together, and tweak a couple of
# small wording/typography points (em dash -> en dash, "usted mismo" ->
# "usted mismo/a").

$d = $word.ActiveDocument

$enDash = [char]0x2013
$emDash = [char]0x2014

# "... que mejor te describe como persona" -> "... que mejor le describa como persona"
$d.Content.Find.Execute("te describe", $true, $false, $false, $false, $false, $true, 1, $false, "le describa", 2) | Out-Null

# "Personas que te conocen" -> "Personas que le conocen"
$d.Content.Find.Execute("Personas que te conocen", $true, $false, $false, $false, $false, $true, 1, $false, "Personas que le conocen", 2) | Out-Null

# "sobre ti —" -> "sobre usted –"
$d.Content.Find.Execute(("sobre ti " + $emDash), $true, $false, $false, $false, $false, $true, 1, $false, ("sobre usted " + $enDash), 2) | Out-Null

# "nosotros queremos saber que piensas sobre" -> "... qué piensa sobre"
$d.Content.Find.Execute("que piensas", $true, $false, $false, $false, $false, $true, 1, $false, "qué piensa", 2) | Out-Null

# "usted mismo" -> "usted mismo/a" (keep this within the underlined run)
$d.Content.Find.Execute("mismo", $true, $false, $false, $false, $false, $true, 1, $false, "mismo/a", 2) | Out-Null

# "Trata de responder ..." -> "Trate de responder ..."
$d.Content.Find.Execute("Trata", $true, $false, $false, $false, $false, $true, 1, $false, "Trate", 2) | Out-Null
